$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (columns reordered: URL/Job Ref moved up, Score(%)/Gap moved up) ---
$ws.Range("A1").Value = "No."
$ws.Range("B1").Value = "Filename"
$ws.Range("C1").Value = "Skills"
$ws.Range("D1").Value = "URL"
$ws.Range("E1").Value = "Job Ref"
$ws.Range("F1").Value = "Company"
$ws.Range("G1").Value = "Position"
$ws.Range("H1").Value = "Score (%)"
$ws.Range("I1").Value = "Gap"
$ws.Range("J1").Value = "Emp Type"
$ws.Range("K1").Value = "Seniority"
$ws.Range("L1").Value = "Industry"
$ws.Range("M1").Value = "Address"
$ws.Range("N1").Value = "Years"

# --- Data rows: newly scanned resumes inserted, renumbered 1-10 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "132.pdf"
$ws.Range("C2").Value = ""

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2.pdf"
$ws.Range("C3").Value = "AutoCAD,Excel,Excellent,MS Office"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "26.pdf"
$ws.Range("C4").Value = "Administration,CRM,Customer Satisfaction,Excel,Excellent"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "3.pdf"
$ws.Range("C5").Value = "Administration,CRM,Customer Satisfaction,Excel,Excellent"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "56.pdf"
$ws.Range("C6").Value = "Communication"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "92.pdf"
$ws.Range("C7").Value = "Excel,Excellent,PowerPoint"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Achyuth Resume_8.docx"
$ws.Range("C8").Value = "Agile,Docker,Excel,Excellent,Java,JavaScript,Software Development,SQL,Web Service,Windows"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Adelina_Erimia_PMP1.docx"
$ws.Range("C9").Value = "Agile,Communication,Construction,Excel,Excellent,ITIL,Leadership,PMP,PowerPoint,project management,Recruiting"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Anil Krishna Mogalaturthi.docx"
$ws.Range("C10").Value = "Administration,Agile,Communication,compliance,Docker,Excel,Excellent,Java,JavaScript,Operating Systems,SQL,Troubleshooting,Web Service,Windows"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "B Shaker-Sr BSA-Scrum Master .docx"
$ws.Range("C11").Value = "Agile,compliance,CRM,Excel,Java,JavaScript,MS Office,Operating Systems,SQL,Web Service,Windows"
